$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target-cluster label change: "Resolving-Mac" -> "Neutrophils" (A4)
$ws.Range("A4").Value = "Neutrophils"

# Row 2 updated values
$ws.Range("I2").Value = 0.09367783019478619
$ws.Range("J2").Value = 0.1170568785762142
$ws.Range("S2").Value = 0.09367783019478619
$ws.Range("T2").Value = 0.1170568785762142

# Row 3 updated values
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8947400000000001
$ws.Range("H3").Value = 1.78948
$ws.Range("I3").Value = 0.5991714967746947
$ws.Range("J3").Value = 0.4991372736284659
$ws.Range("Q3").Value = 0.2493169150266667
$ws.Range("R3").Value = 1.49590149016
$ws.Range("S3").Value = 0.5991714967746947
$ws.Range("T3").Value = 0.4991372736284659

# Row 4 updated values
$ws.Range("G4").Value = 0.4586666666666666
$ws.Range("H4").Value = 1.376
$ws.Range("I4").Value = 0.307150673030519
$ws.Range("J4").Value = 0.3838058477953199
$ws.Range("Q4").Value = 0.1278062435555556
$ws.Range("R4").Value = 1.150256192
$ws.Range("S4").Value = 0.307150673030519
$ws.Range("T4").Value = 0.3838058477953199
